$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that was bumped by one day
# (45204 -> 45205) for every data row, from row 2 through row 261.
$ws.Range("C2:C261").Value = 45205
